$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 15:52"
$ws.Range("A19").Value = "Malaga"
$ws.Range("B19").Value = 2203
$ws.Range("C19").Value = 558
$ws.Range("D19").Value = 1462
$ws.Range("E19").Value = 183
$ws.Range("A20").Value = "Salamanca"
$ws.Range("B20").Value = 2193
$ws.Range("C20").Value = 607
$ws.Range("D20").Value = 1334
$ws.Range("E20").Value = 252
$ws.Range("A21").Value = "Sevilla"
$ws.Range("B21").Value = 2034
$ws.Range("C21").Value = 247
$ws.Range("D21").Value = 1618
$ws.Range("E21").Value = 169
$ws.Range("A22").Value = "A Coruña"
$ws.Range("B22").Value = 1969
$ws.Range("C22").Value = 333
$ws.Range("D22").Value = 1788
$ws.Range("E22").Value = 67
$ws.Range("A23").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B23").Value = 1955
$ws.Range("C23").Value = 4867
$ws.Range("D23").Value = 5101
$ws.Range("E23").Value = 130
$ws.Range("A26").Value = "Granada"
$ws.Range("B26").Value = 1772
$ws.Range("C26").Value = 317
$ws.Range("D26").Value = 1290
$ws.Range("E26").Value = 165
$ws.Range("A27").Value = "Cantabria"
$ws.Range("B27").Value = 1752
$ws.Range("C27").Value = 297
$ws.Range("D27").Value = 1345
$ws.Range("E27").Value = 110
$ws.Range("B32").Value = 1290
$ws.Range("C32").Value = 301
$ws.Range("D32").Value = 926
$ws.Range("B34").Value = 1153
$ws.Range("C34").Value = 176
$ws.Range("D34").Value = 921
$ws.Range("E34").Value = 56
$ws.Range("B35").Value = 1136
$ws.Range("C35").Value = 200
$ws.Range("D35").Value = 831
$ws.Range("E35").Value = 105
$ws.Range("B39").Value = 971
$ws.Range("C39").Value = 194
$ws.Range("D39").Value = 718
$ws.Range("E39").Value = 59
$ws.Range("B48").Value = 481
$ws.Range("C48").Value = 120
$ws.Range("D48").Value = 334
$ws.Range("B51").Value = 409
$ws.Range("C51").Value = 85
$ws.Range("E51").Value = 36
$ws.Range("B52").Value = 328
$ws.Range("C52").Value = 61
$ws.Range("D52").Value = 241
$ws.Range("E52").Value = 26
$ws.Range("B56").Value = 73
$ws.Range("D56").Value = 59
$ws.Range("B57").Value = 68
$ws.Range("D57").Value = 54
$ws.Range("C63").Value = 3
$ws.Range("D63").Value = 4
